$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I0, IF) - copy formatting (bold, border, alignment) from the
# existing header cell H1 so the new columns match the rest of the header row.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I and J
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 6
